$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G: mirror column F (Execute result) into the new Revise column ---
$ws.Range("G2").Value = $ws.Range("F2").Value2
$ws.Range("G3").Value = $ws.Range("F3").Value2
$ws.Range("G4").Value = $ws.Range("F4").Value2
$ws.Range("G5").Value = $ws.Range("F5").Value2
$ws.Range("G6").Value = $ws.Range("F6").Value2

# --- Fix up the running Test Case numbers in column B (5,6 -> 4,5) ---
$ws.Range("B5").Value = 4
$ws.Range("B6").Value = 5

# --- Summary block #1 (Execute / column F) in I3:K6 ---
$ws.Range("I3:K3").Merge()
$ws.Range("I3").Value = "Result"
$ws.Range("I3:K3").Style = "60% - Accent2"
$ws.Range("I3:K3").HorizontalAlignment = -4108
$ws.Range("I3:K3").VerticalAlignment = -4108

$ws.Range("I4:J4").Merge()
$ws.Range("I4").Value = "สรุปผลการทดสอบ"
$ws.Range("K4").Value = "คิดเป็น %"
$ws.Range("I4:K4").Style = "40% - Accent2"
$ws.Range("I4:K4").HorizontalAlignment = -4108
$ws.Range("I4:K4").VerticalAlignment = -4108

$ws.Range("I5").Value = "Pass"
$ws.Range("J5").Formula = '=COUNTIF(F:F,"Pass")'
$ws.Range("K5").Formula = '=TEXT(J5/5,"0.00%")'
$ws.Range("I5:K5").Style = "Good"

$ws.Range("I6").Value = "Fail"
$ws.Range("J6").Formula = '=COUNTIF(F:F,"Fail")'
$ws.Range("K6").Formula = '=TEXT(J6/5,"0.00%")'
$ws.Range("I6:K6").Style = "Bad"

# --- Summary block #2 (Revise / column G) in I9:K12 ---
$ws.Range("I9:K9").Merge()
$ws.Range("I9").Value = "Revise"
$ws.Range("I9:K9").Style = "60% - Accent2"
$ws.Range("I9:K9").HorizontalAlignment = -4108
$ws.Range("I9:K9").VerticalAlignment = -4108

$ws.Range("I10:J10").Merge()
$ws.Range("I10").Value = "สรุปผลการทดสอบ"
$ws.Range("K10").Value = "คิดเป็น %"
$ws.Range("I10:K10").Style = "40% - Accent2"
$ws.Range("I10:K10").HorizontalAlignment = -4108
$ws.Range("I10:K10").VerticalAlignment = -4108

$ws.Range("I11").Value = "Pass"
$ws.Range("J11").Formula = '=COUNTIF(G:G,"Pass")'
$ws.Range("K11").Formula = '=TEXT(J11/5,"0.00%")'
$ws.Range("I11:K11").Style = "Good"

$ws.Range("I12").Value = "Fail"
$ws.Range("J12").Formula = '=COUNTIF(G:G,"Fail")'
$ws.Range("K12").Formula = '=TEXT(J12/5,"0.00%")'
$ws.Range("I12:K12").Style = "Bad"

# --- Selection / window bookkeeping to mirror the saved file state ---
$ws.Range("J13").Select()
